$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab query (cell B2): append a new "Cohort" column to the RETURN clause ---
$newCasesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['West Highland White Terrier'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# --- FilesTab query (cell B4): drop the trailing Study Code line from the RETURN clause ---
$newFilesTabQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['West Highland White Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis  
'@

# --- StatQuery (cells C2, C3, C4 - shared string): replaced entirely with a simpler summary query ---
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['West Highland White Terrier'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("B2").Value2 = $newCasesTabQuery
$ws.Range("B4").Value2 = $newFilesTabQuery
$ws.Range("C2:C4").Value2 = $newStatQuery

# --- Row heights: re-fit to the (smaller) rewritten text ---
$ws.Rows(2).RowHeight = 270
$ws.Rows(3).RowHeight = 225
$ws.Rows(4).RowHeight = 210

# --- Column widths: minor re-fit (character widths) ---
$ws.Columns(1).ColumnWidth = 10.86
$ws.Columns(2).ColumnWidth = 91.43
$ws.Columns(3).ColumnWidth = 74.71
$ws.Columns(4).ColumnWidth = 69.29
$ws.Columns(5).ColumnWidth = 27.71

# --- Selection: user ended up with C2 selected (view scrolled back to top) ---
$ws.Range("C2").Select()
